$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.959.26"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.449.91"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'523.47"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'131.14"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "2.454.84"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("D10").Value = "'0.0982"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "'4.95"
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").Value = "'0.324"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "2.884.21"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "57.858.72"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "'21.73"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "2.448.48"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "'10.28"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "'314.52"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").Value = "'6.14"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'64.91"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Value = "'0.402"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").Value = "'7.22"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").Value = "'173.11"
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "'6.20"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "'17.80"
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("D38").Value = "'3.79"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").Value = "'36.24"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "'1.45"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").Value = "'0.797"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "'263.06"
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("D44").Value = "'0.587"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "'4.81"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'122.69"
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0920"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'0.0497"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "'0.0212"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'16.96"
$ws.Range("E50").Value = "  -4.31%  "
$ws.Range("D51").Value = "'16.25"
$ws.Range("E51").Value = "  -3.77%  "
